$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "Förändrad" (Changed) date column for rows 2-18 from 2023-09-05 (45174) to 2023-09-06 (45175)
$ws.Range("C2:C18").Value = 45175
